$wb = $excel.ActiveWorkbook

# Sheet: general
$ws = $wb.Worksheets.Item("general")
$ws.Range("B3").Value = 272.5632163890596
$ws.Range("B4").Value = 0.01800012588500977
$ws.Range("B6").Value = 40.67321638905945
$ws.Range("B7").Value = 3.515920931989228
$ws.Range("B8").Value = 3.515920931989228
$ws.Range("B9").Value = 231.89

# Sheet: alpha
$ws = $wb.Worksheets.Item("alpha")
$ws.Range("A2").Value = 4
$ws.Range("A3").Value = 4
$ws.Range("A4").Value = 4
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 1
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 1
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = 1

# Sheet: x
$ws = $wb.Worksheets.Item("x")
$ws.Range("B2").Value = 3
$ws.Range("B4").Value = 9
$ws.Range("B6").Value = 11
$ws.Range("B7").Value = 6
$ws.Range("B8").Value = 12
$ws.Range("B9").Value = 10
$ws.Range("B10").Value = 2
$ws.Range("B12").Value = 7
$ws.Range("B14").Value = 4

# Sheet: TBar
$ws = $wb.Worksheets.Item("TBar")
$ws.Range("B3").Value = 35.46541289070538
$ws.Range("B4").Value = 34.69770569366315
$ws.Range("B5").Value = 30
$ws.Range("B6").Value = 38.43347677669323
$ws.Range("B7").Value = 31.94987179065701
$ws.Range("B8").Value = 30.34885527085025
$ws.Range("B9").Value = 34.86507964991324
$ws.Range("B11").Value = 30
$ws.Range("B12").Value = 34.87971820301381
$ws.Range("B13").Value = 37.27819014430416
$ws.Range("B14").Value = 41.94939770868245
$ws.Range("B15").Value = 40.80207851635592

# Sheet: y
$ws = $wb.Worksheets.Item("y")
$ws.Range("A2").Value = 4
$ws.Range("B2").Value = 12
$ws.Range("A3").Value = 4
$ws.Range("B3").Value = 12
$ws.Range("A4").Value = 4
$ws.Range("B4").Value = 12
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 1
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 12
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = 1
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 12
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = 1

# Sheet: Q
$ws = $wb.Worksheets.Item("Q")
$ws.Range("C7").Value = 260.9450000000011
$ws.Range("C8").Value = 281.47
$ws.Range("C9").Value = 251.0650000000011
$ws.Range("C10").Value = 272.6950000000011
$ws.Range("C11").Value = 256.5900000000011
$ws.Range("C12").Value = 235.775
$ws.Range("C13").Value = 229.025
$ws.Range("C14").Value = 213.42
$ws.Range("C15").Value = 226.76
$ws.Range("C16").Value = 221.56
$ws.Range("C17").Value = 46.91999999999942
$ws.Range("C18").Value = 36.10499999999942
$ws.Range("C19").Value = 34.91499999999942
$ws.Range("C20").Value = 37.48999999999942
$ws.Range("C21").Value = 39.43499999999941
$ws.Range("C22").Value = 285.945
$ws.Range("C23").Value = 303.02
$ws.Range("C24").Value = 296.175
$ws.Range("C25").Value = 307.985
$ws.Range("C26").Value = 291.35
$ws.Range("C27").Value = 155.3650000000007
$ws.Range("C28").Value = 164.2850000000007
$ws.Range("C29").Value = 152.1250000000008
$ws.Range("C30").Value = 163.3750000000007
$ws.Range("C31").Value = 157.6200000000007
$ws.Range("C32").Value = 154.3
$ws.Range("C33").Value = 148.3449999999993
$ws.Range("C34").Value = 128.7049999999993
$ws.Range("C35").Value = 146.3249999999992
$ws.Range("C36").Value = 134.2149999999993
$ws.Range("C37").Value = 151
$ws.Range("C38").Value = 163.8800000000003
$ws.Range("C39").Value = 137.7950000000003
$ws.Range("C40").Value = 153.5400000000003
$ws.Range("C41").Value = 140.7850000000003
$ws.Range("C47").Value = 85.48500000000051
$ws.Range("C48").Value = 87.9650000000005
$ws.Range("C49").Value = 79.71500000000052
$ws.Range("C50").Value = 90.7300000000005
$ws.Range("C51").Value = 84.73000000000052
$ws.Range("C52").Value = 213.315
$ws.Range("C53").Value = 222.965
$ws.Range("C54").Value = 213.855
$ws.Range("C55").Value = 224.03
$ws.Range("C56").Value = 210.535
$ws.Range("C57").Value = 235.775
$ws.Range("C58").Value = 229.025
$ws.Range("C59").Value = 213.42
$ws.Range("C60").Value = 226.76
$ws.Range("C61").Value = 221.56
$ws.Range("C62").Value = 285.945
$ws.Range("C63").Value = 303.02
$ws.Range("C64").Value = 296.175
$ws.Range("C65").Value = 307.985
$ws.Range("C66").Value = 291.35
$ws.Range("C67").Value = 260.9450000000011
$ws.Range("C68").Value = 281.47
$ws.Range("C69").Value = 251.0650000000011
$ws.Range("C70").Value = 272.6950000000011
$ws.Range("C71").Value = 256.5900000000011

# Sheet: R
$ws = $wb.Worksheets.Item("R")
$ws.Range("C2").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("C7").Value = 10.945
$ws.Range("C8").Value = 28.02
$ws.Range("C9").Value = 21.175
$ws.Range("C10").Value = 32.985
$ws.Range("C11").Value = 16.35
$ws.Range("C13").Value = 6.470000000000027

# Sheet: rho
$ws = $wb.Worksheets.Item("rho")
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 2
$ws.Range("A3").Value = 4
$ws.Range("B3").Value = 1
$ws.Range("A4").Value = 4
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 1
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = 1
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 4
$ws.Range("C6").Value = 1
$ws.Range("A7").Value = 4
$ws.Range("B7").Value = 5
$ws.Range("C7").Value = 1
